$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","J","K","M","N","O")

$data = @{
    2 = @(4.915629811074257,8.583274659809444,13.05105998099303,34.31405872449314,3.661337967947568,9.892769714158318,14.51406234845141,17.38098029013006,18.35145926791868,26.0842163145504)
    3 = @(4.750894200201358,8.572011588633375,13.06757547913093,34.35245704808164,3.663666144554797,9.917952408448256,14.13593058514326,17.23843288551542,18.40513111931306,26.1437608241578)
    4 = @(4.648291909041469,8.566166213489042,13.07973477895648,34.38497206122835,3.665171642017205,9.934520280079946,13.9006312070146,17.15309096936968,18.43993193558751,26.18673831074121)
    5 = @(4.606190335121231,8.564055142500637,13.08519753520319,34.40046493900022,3.665804314275285,9.941550167235782,13.80410058780824,17.1188932248781,18.45457881296792,26.20586066810353)
    6 = @(4.599184177112573,8.563721018105047,13.0861352899141,34.40317281391251,3.66591052865909,9.9427342958725,13.78803712479076,17.11325061615736,18.45703904765439,26.20913293038522)
    7 = @(4.647725184292873,8.566136643353344,13.07980639566985,34.38517193002731,3.665180096756405,9.9346139601445,13.89933177723451,17.15262737905161,18.44012758327679,26.18698969485867)
    8 = @(4.859174163822052,8.579170185937409,13.05633563592205,34.3254407014453,3.662124989462022,9.901223441520134,14.38441370047176,17.33139574324345,18.36958279731661,26.10341259727855)
    9 = @(5.259272899247185,8.61313416879001,13.02631959348161,34.27940411879399,3.656734027361218,9.844502188337106,15.30503926839248,17.69787953951602,18.24584507921249,25.99063699973582)
    10 = @(5.5407350462609,8.643089094656013,13.01401326793241,34.28909916195072,3.653135155497707,9.808147346495071,15.95555556709328,17.97491704665595,18.16377002586783,25.93919851941967)
    11 = @(5.665420649230609,8.657773590814601,13.01052633362513,34.30297352153691,3.65157566791075,9.792759129736044,16.24453000137759,18.10221916390116,18.12833676272603,25.92266497945402)
    12 = @(5.712107463491328,8.663483660778885,13.00950891118169,34.30958736990895,3.650996233478555,9.787097038635793,16.35285793407395,18.15057312110622,18.11519178534389,25.91739433402049)
    13 = @(5.702076834964631,8.662247296910888,13.00971456576707,34.30810250491714,3.651120531854119,9.788309132570157,16.3295780335355,18.14015325017656,18.11801067019107,25.91848537956207)
    14 = @(5.669272450896946,8.658240384310373,13.0104365608145,34.30349040184576,3.651527775164464,9.79228999877658,16.25346480755405,18.10619452563738,18.12724985427546,25.92221149945316)
    15 = @(5.649108640352314,8.655805403256345,13.01091824416176,34.30084240818479,3.651778668563029,9.794749887022812,16.206697101154,18.08541195124654,18.13294462760649,25.92462289231453)
    16 = @(5.532514959051118,8.642150483067681,13.01428358318044,34.28838282420153,3.653238629429687,9.809176118189585,15.93652172774087,17.96662014909271,18.16612389632674,25.94041744933344)
    17 = @(5.460094757809501,8.634042659475771,13.016888469152,34.28316292721982,3.654154117082712,9.818320455070632,15.76892331110348,17.89404643387112,18.18696513050941,25.95186770654779)
    18 = @(5.41812742322755,8.629478996519087,13.01858549632151,34.28105151702578,3.654687994236929,9.823688273134383,15.67187656633417,17.85242625785801,18.19913163315733,25.95909971844767)
    19 = @(5.403865733978002,8.627951031126045,13.01919423509367,34.28048965902119,3.654870013615672,9.825524321510587,15.63891004407627,17.83835648818072,18.20328180131558,25.96165922898252)
    20 = @(5.467836777470028,8.634895447634459,13.01659060706007,34.28362638878492,3.654055905421019,9.817335824425937,15.78683233444321,17.90175963331815,18.18472800759893,25.95058191260584)
    21 = @(5.678922594829841,8.659413280354407,13.01021627521169,34.30480819419359,3.651407856873572,9.791116243234681,16.27585173433452,18.11616531339613,18.12452868391611,25.92109015210379)
    22 = @(5.813778412738436,8.676306378464794,13.00781606605836,34.32657710881706,3.649741933428317,9.774942407180703,16.58899940869351,18.25713432306566,18.08677484163246,25.90758847461564)
    23 = @(5.742101161869488,8.667211616245426,13.00893575822422,34.31423410153672,3.650625164221111,9.783486728339977,16.42248835930276,18.18183134184487,18.10677957063532,25.91426552939728)
    24 = @(5.464337640898624,8.634509597798409,13.01672464922628,34.28341408671974,3.654100283387007,9.817780631354452,15.77873781317928,17.89827217093467,18.18573883576302,25.95116119829437)
    25 = @(5.152990318160145,8.603058263543041,13.03272676011969,34.28422423577511,3.65812859275275,9.85891138442352,15.06004142618431,17.69787953951602,18.27776322510569,26.0156450896837)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowVals[$i]
    }
}
